$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of portfolio data appended below the last existing row (row 54).
# Column A holds the date as text (matches the existing rows, which store
# dates as plain text rather than Excel date serials), so force a text
# number format before assigning the value to stop Excel's COM layer from
# auto-converting the "yyyy-mm-dd" looking string into a date serial.
$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value = "2025-10-09"
$ws.Range("A55").ClearFormats()
$ws.Range("B55").Value = 53.16999816894531
$ws.Range("C55").Value = 681.0999755859375
$ws.Range("D55").Value = 345.5
